$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Rows 35 and 36 got their match records (B:AB) swapped between each other.
# Column A (the row index) stays put on each row.
# ---------------------------------------------------------------------------
$ws.Range("B35").Value = 6847042
$ws.Range("B36").Value = 6847044
$ws.Range("C35").Value = "Austria Bundesliga"
$ws.Range("C36").Value = "Austria Bundesliga"
$ws.Range("D35").Value = 45165.5
$ws.Range("D36").Value = 45165.5
$ws.Range("E35").Value = "LASK Linz"
$ws.Range("E36").Value = "Hartberg"
$ws.Range("F35").Value = "FK Austria Vienna"
$ws.Range("F36").Value = "Austria Klagenfurt"
$ws.Range("G35").Value = 2
$ws.Range("G36").Value = 0
$ws.Range("H35").Value = 0
$ws.Range("H36").Value = 3
$ws.Range("I35").Value = "H"
$ws.Range("I36").Value = "A"
$ws.Range("J35").Value = 1.909
$ws.Range("J36").Value = 2.75
$ws.Range("K35").Value = 3.5
$ws.Range("K36").Value = 3.25
$ws.Range("L35").Value = 3.25
$ws.Range("L36").Value = 2.25
$ws.Range("M35").Value = 1.85
$ws.Range("M36").Value = 2.1
$ws.Range("N35").Value = 3.8
$ws.Range("N36").Value = 3.5
$ws.Range("O35").Value = 4
$ws.Range("O36").Value = 3.3
$ws.Range("P35").Value = -0.5
$ws.Range("P36").Value = -0.25
$ws.Range("Q35").Value = 1.85
$ws.Range("Q36").Value = 1.8
$ws.Range("R35").Value = 2
$ws.Range("R36").Value = 2.05
$ws.Range("S35").Value = 3
$ws.Range("S36").Value = 2.5
$ws.Range("T35").Value = 1.9
$ws.Range("T36").Value = 1.875
$ws.Range("U35").Value = 1.95
$ws.Range("U36").Value = 1.975
$ws.Range("V35").Value = 0.8500000000000001
$ws.Range("V36").Value = -1
$ws.Range("W35").Value = -1
$ws.Range("W36").Value = -1
$ws.Range("X35").Value = -1
$ws.Range("X36").Value = 2.3
$ws.Range("Y35").Value = 0.8500000000000001
$ws.Range("Y36").Value = -1
$ws.Range("Z35").Value = -1
$ws.Range("Z36").Value = 1.05
$ws.Range("AA35").Value = -1
$ws.Range("AA36").Value = 0.875
$ws.Range("AB35").Value = 0.95
$ws.Range("AB36").Value = -1

# ---------------------------------------------------------------------------
# Rows 53 and 54 likewise had their match records (B:AB) swapped.
# ---------------------------------------------------------------------------
$ws.Range("B53").Value = 6847056
$ws.Range("B54").Value = 6847053
$ws.Range("C53").Value = "Austria Bundesliga"
$ws.Range("C54").Value = "Austria Bundesliga"
$ws.Range("D53").Value = 45193.39583333334
$ws.Range("D54").Value = 45193.39583333334
$ws.Range("E53").Value = "SCR Altach"
$ws.Range("E54").Value = "LASK Linz"
$ws.Range("F53").Value = "FK Austria Vienna"
$ws.Range("F54").Value = "Hartberg"
$ws.Range("G53").Value = 2
$ws.Range("G54").Value = 0
$ws.Range("H53").Value = 1
$ws.Range("H54").Value = 0
$ws.Range("I53").Value = "H"
$ws.Range("I54").Value = "D"
$ws.Range("J53").Value = 3.75
$ws.Range("J54").Value = 1.571
$ws.Range("K53").Value = 3.5
$ws.Range("K54").Value = 4.333
$ws.Range("L53").Value = 1.95
$ws.Range("L54").Value = 5
$ws.Range("M53").Value = 3
$ws.Range("M54").Value = 1.727
$ws.Range("N53").Value = 3.3
$ws.Range("N54").Value = 4
$ws.Range("O53").Value = 2.4
$ws.Range("O54").Value = 4.5
$ws.Range("P53").Value = 0
$ws.Range("P54").Value = -0.75
$ws.Range("Q53").Value = 2.1
$ws.Range("Q54").Value = 1.95
$ws.Range("R53").Value = 1.775
$ws.Range("R54").Value = 1.9
$ws.Range("S53").Value = 2.5
$ws.Range("S54").Value = 2.75
$ws.Range("T53").Value = 1.925
$ws.Range("T54").Value = 1.8
$ws.Range("U53").Value = 1.925
$ws.Range("U54").Value = 2.05
$ws.Range("V53").Value = 2
$ws.Range("V54").Value = -1
$ws.Range("W53").Value = -1
$ws.Range("W54").Value = 3
$ws.Range("X53").Value = -1
$ws.Range("X54").Value = -1
$ws.Range("Y53").Value = 1.1
$ws.Range("Y54").Value = -1
$ws.Range("Z53").Value = -1
$ws.Range("Z54").Value = 0.8999999999999999
$ws.Range("AA53").Value = 0.925
$ws.Range("AA54").Value = -1
$ws.Range("AB53").Value = -1
$ws.Range("AB54").Value = 1.05

# ---------------------------------------------------------------------------
# Rows 166 and 167 likewise had their match records (B:AB) swapped.
# ---------------------------------------------------------------------------
$ws.Range("B166").Value = 7948274
$ws.Range("B167").Value = 7948276
$ws.Range("C166").Value = "Austria Bundesliga"
$ws.Range("C167").Value = "Austria Bundesliga"
$ws.Range("D166").Value = 45402.5
$ws.Range("D167").Value = 45402.5
$ws.Range("E166").Value = "Wolfsberger AC"
$ws.Range("E167").Value = "SCR Altach"
$ws.Range("F166").Value = "Austria Lustenau"
$ws.Range("F167").Value = "FC Blau Weiss Linz"
$ws.Range("G166").Value = 1
$ws.Range("G167").Value = 2
$ws.Range("H166").Value = 1
$ws.Range("H167").Value = 2
$ws.Range("I166").Value = "D"
$ws.Range("I167").Value = "D"
$ws.Range("J166").Value = 1.533
$ws.Range("J167").Value = 2.1
$ws.Range("K166").Value = 4
$ws.Range("K167").Value = 3.2
$ws.Range("L166").Value = 6
$ws.Range("L167").Value = 3.6
$ws.Range("M166").Value = 1.666
$ws.Range("M167").Value = 2
$ws.Range("N166").Value = 3.8
$ws.Range("N167").Value = 3
$ws.Range("O166").Value = 5
$ws.Range("O167").Value = 4.75
$ws.Range("P166").Value = -0.75
$ws.Range("P167").Value = -0.5
$ws.Range("Q166").Value = 1.925
$ws.Range("Q167").Value = 2
$ws.Range("R166").Value = 1.925
$ws.Range("R167").Value = 1.85
$ws.Range("S166").Value = 2.25
$ws.Range("S167").Value = 1.75
$ws.Range("T166").Value = 1.95
$ws.Range("T167").Value = 1.775
$ws.Range("U166").Value = 1.9
$ws.Range("U167").Value = 2.1
$ws.Range("V166").Value = -1
$ws.Range("V167").Value = -1
$ws.Range("W166").Value = 2.8
$ws.Range("W167").Value = 2
$ws.Range("X166").Value = -1
$ws.Range("X167").Value = -1
$ws.Range("Y166").Value = -1
$ws.Range("Y167").Value = -1
$ws.Range("Z166").Value = 0.925
$ws.Range("Z167").Value = 0.8500000000000001
$ws.Range("AA166").Value = -0.5
$ws.Range("AA167").Value = 0.7749999999999999
$ws.Range("AB166").Value = 0.45
$ws.Range("AB167").Value = -1

# ---------------------------------------------------------------------------
# Rows 182, 183, 185 and 187 are future (not yet played) fixtures whose
# closing odds got refreshed with newer prices.
# ---------------------------------------------------------------------------
$ws.Range("M182").Value = 5.5
$ws.Range("O182").Value = 1.666
$ws.Range("Q182").Value = 1.975
$ws.Range("R182").Value = 1.875
$ws.Range("T182").Value = 1.975
$ws.Range("U182").Value = 1.875

$ws.Range("T183").Value = 2.1
$ws.Range("U183").Value = 1.775

$ws.Range("M185").Value = 1.4
$ws.Range("Q185").Value = 1.9
$ws.Range("R185").Value = 1.95
$ws.Range("T185").Value = 1.85
$ws.Range("U185").Value = 2

$ws.Range("M187").Value = 4.5
$ws.Range("N187").Value = 3.8
$ws.Range("O187").Value = 1.7
$ws.Range("P187").Value = 0.75
$ws.Range("Q187").Value = 1.925
$ws.Range("R187").Value = 1.925
